$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "Total" row at row 60 summing columns B and C (rows 2-59)
$ws.Range("A60").Value = "Total"
$ws.Range("B60").Value = 158
$ws.Range("C60").Value = 111
